$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; D="39.884.56"; E="  -4.35%  "}
    @{Row=3; D="2.327.16"; E="  -5.72%  "}
    @{Row=4; D="0.999"; E="  -0.10%  "}
    @{Row=5; D="308.16"; E="  -3.88%  "}
    @{Row=6; D="84.58"; E="  -8.23%  "}
    @{Row=7; D="0.531"; E="  -3.30%  "}
    @{Row=8; D=$null; E="  +0.03%  "}
    @{Row=9; D="0.484"; E="  -4.62%  "}
    @{Row=10; D="0.0812"; E="  -4.88%  "}
    @{Row=11; D="29.99"; E="  -8.80%  "}
    @{Row=12; D=$null; E="  +0.49%  "}
    @{Row=13; D="2.685.55"; E="  -5.76%  "}
    @{Row=14; D="6.39"; E="  -6.99%  "}
    @{Row=15; D="14.65"; E="  -5.42%  "}
    @{Row=16; D="2.327.64"; E="  -6.09%  "}
    @{Row=17; D="0.753"; E="  -4.58%  "}
    @{Row=18; D="39.861.50"; E="  -4.22%  "}
    @{Row=19; D="0.0₃0902"; E="  -3.92%  "}
    @{Row=20; D="6.06"; E="  -5.73%  "}
    @{Row=21; D="67.47"; E="  -6.08%  "}
    @{Row=22; D="10.59"; E="  -5.41%  "}
    @{Row=23; D="234.77"; E="  -1.90%  "}
    @{Row=24; D="2.55"; E="  -7.26%  "}
    @{Row=25; D=$null; E="  +0.20%  "}
    @{Row=26; D=$null; E="  -6.78%  "}
    @{Row=27; D="23.25"; E="  -6.09%  "}
    @{Row=28; D=$null; E="  -4.28%  "}
    @{Row=29; D="9.25"; E="  -4.61%  "}
    @{Row=30; D="35.04"; E="  -2.84%  "}
    @{Row=31; D="151.75"; E="  -2.31%  "}
    @{Row=32; D=$null; E="  -0.06%  "}
    @{Row=33; D="5.09"; E="  -6.13%  "}
    @{Row=34; D=$null; E="  -4.43%  "}
    @{Row=35; D="0.0719"; E="  -5.69%  "}
    @{Row=36; D=$null; E="  -2.53%  "}
    @{Row=37; D="0.0995"; E="  -3.40%  "}
    @{Row=38; D="2.74"; E="  -5.75%  "}
    @{Row=39; D="15.56"; E="  -8.41%  "}
    @{Row=40; D="1.70"; E="  -6.98%  "}
    @{Row=41; D="3.82"; E="  -4.01%  "}
    @{Row=42; D=$null; E="  -3.10%  "}
    @{Row=43; D="1.946.06"; E="  -2.78%  "}
    @{Row=44; D="0.0265"; E="  -5.82%  "}
    @{Row=45; D="17.54"; E="  -5.52%  "}
    @{Row=46; D="9.34"; E="  -1.18%  "}
    @{Row=47; D=$null; E="  -9.46%  "}
    @{Row=48; D="2.553.49"; E="  -6.54%  "}
    @{Row=49; D="92.60"; E="  -4.60%  "}
    @{Row=50; D="70.41"; E=$null}
    @{Row=51; D="50.15"; E="  -3.81%  "}
)
foreach ($item in $updates) {
    $r = $item.Row
    if ($null -ne $item.D) {
        $dCell = $ws.Range("D$r")
        $dCell.NumberFormat = "@"
        $dCell.Value = $item.D
        $dCell.Style = "Normal"
    }
    if ($null -ne $item.E) {
        $ws.Range("E$r").Value = $item.E
    }
}
